# Applies the authored edit to ProjectSubmissionTemplate-210521-090339.pptx:
#   1. Bumps the "datetimeFigureOut" date placeholder text from 6/27/2021 to
#      6/28/2021 on the slide master and all five slide layouts.
#   2. On slide 14's SWOT grid:
#        - Shape "Rectangle 10" ("An innovative approach ...") loses the
#          trailing "(online, advanced software..)" qualifier and the
#          remaining text is split into two runs ("...in Haiti " / "for
#          Decision-making "); the textbox autosizes shorter.
#        - Shape "Rectangle 11" loses its second bullet ("Lack of
#          resources...") while keeping the empty paragraph; the textbox
#          autosizes shorter.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the date placeholder everywhere it appears: slide master +
#    all custom (slide) layouts.
# ---------------------------------------------------------------------
function Update-DateHolder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -eq "6/27/2021") {
                $shp.TextFrame.TextRange.Text = "6/28/2021"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateHolder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateHolder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 14 SWOT grid edits.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(14)

# --- Shape "Rectangle 10": shorten the title & split into two runs ---
$titleShape = $slide.Shapes.Item("Rectangle 10")
$titleRange = $titleShape.TextFrame.TextRange

$fullText = $titleRange.Text
$newRun1Text = "An innovative approach to teaching Data in Haiti "
$newRun2Text = "for Decision-making "

# Remove the trailing qualifier text, keeping only the two sentences above.
$keepLength = $newRun1Text.Length + $newRun2Text.Length
if ($fullText.Length -gt $keepLength) {
    $trailing = $titleRange.Characters($keepLength + 1, $fullText.Length - $keepLength)
    $trailing.Text = ""
}

# Re-assign the second half of the remaining text so it becomes its own run
# (same formatting is inherited from the original run it is carved out of).
$secondHalf = $titleRange.Characters($newRun1Text.Length + 1, $newRun2Text.Length)
$secondHalf.Text = $newRun2Text

# Shrink the textbox to the new autofit height.
$titleShape.Height = 430887 / 12700

# --- Shape "Rectangle 11": drop the second bullet's text ---
$swShape = $slide.Shapes.Item("Rectangle 11")
$swRange = $swShape.TextFrame.TextRange
$secondParagraph = $swRange.Paragraphs(2, 1)
$secondParagraph.Text = ""

# Shrink the textbox to the new autofit height.
$swShape.Height = 584775 / 12700
